$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 3.2
$ws.Range("L2").Value = 1.5
$ws.Range("M2").Value = 2.18
$ws.Range("R2").Value = 55
$ws.Range("S2").Value = 60
$ws.Range("T2").Value = 7
$ws.Range("H8").Value = 2.77
$ws.Range("O8").Value = 11
$ws.Range("Y8").Value = 5.8
$ws.Range("AD8").Value = 29
$ws.Range("G9").Value = 1.7
$ws.Range("H9").Value = 3.5
$ws.Range("I9").Value = 5.25
$ws.Range("K9").Value = 1.62
$ws.Range("O9").Value = 7
$ws.Range("T9").Value = 8
$ws.Range("W9").Value = 67
$ws.Range("Z9").Value = 26
$ws.Range("AA9").Value = 17
$ws.Range("AB9").Value = 51
$ws.Range("AC9").Value = 41
$ws.Range("AF9").Value = 8
$ws.Range("AI9").Value = 2.1
$ws.Range("AJ9").Value = 1.67
$ws.Range("G10").Value = 1.44
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 1.57
$ws.Range("K10").Value = 2.35
$ws.Range("Q10").Value = 10
$ws.Range("R10").Value = 11
$ws.Range("U10").Value = 9.5
$ws.Range("Y10").Value = 19
$ws.Range("Z10").Value = 34
$ws.Range("AA10").Value = 19
$ws.Range("AB10").Value = 67
$ws.Range("AG10").Value = 1.17
$ws.Range("AH10").Value = 5
$ws.Range("G12").Value = 2.65
$ws.Range("I12").Value = 2.6
$ws.Range("L12").Value = 1.5
$ws.Range("M12").Value = 2.25
$ws.Range("N12").Value = 6.5
$ws.Range("O12").Value = 11.5
$ws.Range("P12").Value = 10.75
$ws.Range("Q12").Value = 30
$ws.Range("R12").Value = 28
$ws.Range("S12").Value = 50
$ws.Range("U12").Value = 6.1
$ws.Range("V12").Value = 19
$ws.Range("Y12").Value = 6.4
$ws.Range("Z12").Value = 11.25
$ws.Range("AA12").Value = 10.75
$ws.Range("AB12").Value = 28
$ws.Range("AC12").Value = 27
$ws.Range("AE14").Value = 1.08
$ws.Range("AG14").Value = 1.49
$ws.Range("AH14").Value = 2.45
$ws.Range("J17").Value = 2.15
$ws.Range("K17").Value = 1.67
$ws.Range("G19").Value = 1.4
$ws.Range("G20").Value = 1.38
$ws.Range("H20").Value = 3.9
$ws.Range("I20").Value = 8
$ws.Range("N20").Value = 5
$ws.Range("O20").Value = 5.5
$ws.Range("Q20").Value = 8.5
$ws.Range("U20").Value = 8.5
$ws.Range("V20").Value = 26
$ws.Range("AD20").Value = 81
$ws.Range("AI20").Value = 2.5
$ws.Range("AJ20").Value = 1.5
$ws.Range("G21").Value = 1.45
$ws.Range("I22").Value = 1.53
$ws.Range("G23").Value = 2.35
$ws.Range("I23").Value = 3
$ws.Range("N23").Value = 6
$ws.Range("O23").Value = 10
$ws.Range("R23").Value = 23
$ws.Range("T23").Value = 6
$ws.Range("W23").Value = 81
$ws.Range("AA23").Value = 13
$ws.Range("AC23").Value = 34
$ws.Range("G24").Value = 1.65
$ws.Range("H24").Value = 3.7
$ws.Range("P24").Value = 9
$ws.Range("T24").Value = 9
$ws.Range("N26").Value = 5.7
$ws.Range("P26").Value = 8.75
$ws.Range("T26").Value = 6.8
$ws.Range("U26").Value = 6.1
$ws.Range("Z26").Value = 19.5
$ws.Range("H27").Value = 2.95
$ws.Range("S27").Value = 45
$ws.Range("T27").Value = 6.7

